# "fixed set_invalid_row method" - row 7 on Sheet1 is a perfectly valid row,
# but its B column (Email) was left blank and un-linked by the old buggy
# logic. Give it the same value + mailto hyperlink that every other row in
# column B already has.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")   # the "Users" sheet (also $wb.ActiveSheet)

$emailCell = $ws.Range("B7")

# Restore the missing email value.
$emailCell.Value = "email1@test.com"

# Re-create the hyperlink the same way it's wired up for the other rows.
[void]$ws.Hyperlinks.Add($emailCell, "mailto:email1@test.com", `
    [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, `
    "email1@test.com")

# Match the formatting (blue/underlined font) already used by the other
# hyperlinked cells in column B, rather than leaving Excel's freshly
# auto-generated "Hyperlink" style on the cell.
$refCell = $ws.Range("B2")
$emailCell.Font.Name = $refCell.Font.Name()
$emailCell.Font.Color = $refCell.Font.Color()
$emailCell.Font.Underline = $refCell.Font.Underline()

# The active selection shifted from A8 to B8 while fixing this up.
[void]$ws.Range("B8").Select()
